$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("situation names")

# Increment each value in A2:A36 by 1 (odd -> even), fixing the situation number.
for ($r = 2; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = [double]$cell.Value2
    $cell.Value = $current + 1
}

# Update the selection to match the resulting author selection (B24 single cell).
$ws.Activate()
$ws.Range("B24").Select()
